$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells hold text (e.g. "26.260.17", "0.5113", percentages) rather than
# numbers, so force text format before assigning values to avoid Excel auto-converting
# them into numeric/scientific values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.260.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -6.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.672.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.85%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5113"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -11.71%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2662"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06364"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.61"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07363"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.55%  "

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.556"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.669.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5809"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.893.71"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008584"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -12.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.319.55"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.943"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "189.76"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.208"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.71"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.673"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1176"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.73"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05825"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.320"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.530"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.517"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.652"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.009"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5997"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.45%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.644"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01617"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.010"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.082.83"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8596"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.80%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.816.66"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.87%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.89"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.009"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4293"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05184"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.62%  "
